# Generate Report for Handoff
# This script rewrites the localization-status report: two new files
# (callerMd1.md / callerMd2.md) replace the old png/md sample rows, and a
# fourth data row is appended to every sheet.

$wb = $excel.ActiveWorkbook

function Remove-CellHyperlink {
    param($ws, $addr)
    $target = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $target = $hl
        }
    }
    if ($target -ne $null) {
        $target.Delete()
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

Remove-CellHyperlink $ws1 '$A$2'
Remove-CellHyperlink $ws1 '$A$3'
Remove-CellHyperlink $ws1 '$A$4'

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-23 09:09:56"
$ws1.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-23 09:09:56"
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-23 09:09:56"
$ws1.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-23 09:09:56"
$ws1.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

Remove-CellHyperlink $ws2 '$A$2'
Remove-CellHyperlink $ws2 '$D$2'
Remove-CellHyperlink $ws2 '$A$3'
Remove-CellHyperlink $ws2 '$D$3'
Remove-CellHyperlink $ws2 '$A$4'
Remove-CellHyperlink $ws2 '$D$4'

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca51f9130789d54479a6ba6a76edff8a0430eb36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca51f9130789d54479a6ba6a76edff8a0430eb36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca51f9130789d54479a6ba6a76edff8a0430eb36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca51f9130789d54479a6ba6a76edff8a0430eb36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf")

$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-23 09:09:52"
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J2").Value = "Include"
$ws2.Range("K2").Value = "e2e\callerMd1.md,`ne2e\callerMd2.md"

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "2016-03-23 09:09:52"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J3").Value = "Include"
$ws2.Range("K3").Value = "e2e\callerMd1.md"

$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("E4").Value = "2016-03-23 09:09:52"
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws2.Range("J4").Value = "Include"
$ws2.Range("K4").ClearContents()

$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("E5").Value = "2016-03-23 09:09:52"
$ws2.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I5").Value = "e2e\calleeMd1.md"
$ws2.Range("J5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

Remove-CellHyperlink $ws3 '$A$2'
Remove-CellHyperlink $ws3 '$D$2'
Remove-CellHyperlink $ws3 '$A$3'
Remove-CellHyperlink $ws3 '$D$3'
Remove-CellHyperlink $ws3 '$A$4'
Remove-CellHyperlink $ws3 '$D$4'

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b147b0ea0455bbfa609814bf212f5f707f2c0deb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b147b0ea0455bbfa609814bf212f5f707f2c0deb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b147b0ea0455bbfa609814bf212f5f707f2c0deb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d342627d14ceb694934136a53cedab65d84c28c2/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b147b0ea0455bbfa609814bf212f5f707f2c0deb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf")

$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-23 09:09:56"
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J2").Value = "Include"
$ws3.Range("K2").Value = "e2e\callerMd1.md,`ne2e\callerMd2.md"

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "2016-03-23 09:09:56"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J3").Value = "Include"
$ws3.Range("K3").Value = "e2e\callerMd1.md"

$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("E4").Value = "2016-03-23 09:09:56"
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws3.Range("J4").Value = "Include"
$ws3.Range("K4").ClearContents()

$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("E5").Value = "2016-03-23 09:09:56"
$ws3.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I5").Value = "e2e\calleeMd1.md"
$ws3.Range("J5").Value = "Include"

Write-Host "Report regenerated for handoff."
